$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date from 45205 to 45206 for every existing
# data row (2 through 499).
$ws.Range("C2:C499").Value2 = 45206

# Row 499 becomes consistent with the rest of the data rows and gets an
# explicit (custom) row height of 15.
$ws.Rows.Item(499).RowHeight = 15

# Add the new row 500 with the new entry's data.
$ws.Cells.Item(500, 1).Value2 = "A 48150-2023"

$ws.Cells.Item(500, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(500, 2).Value2 = 45205

$ws.Cells.Item(500, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(500, 3).Value2 = 45206

$ws.Cells.Item(500, 4).Value2 = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(500, 5).Value2 = "BORÅS"

$ws.Cells.Item(500, 7).Value2 = 3.1
$ws.Cells.Item(500, 8).Value2 = 0
$ws.Cells.Item(500, 9).Value2 = 0
$ws.Cells.Item(500, 10).Value2 = 0
$ws.Cells.Item(500, 11).Value2 = 0
$ws.Cells.Item(500, 12).Value2 = 0
$ws.Cells.Item(500, 13).Value2 = 0
$ws.Cells.Item(500, 14).Value2 = 0
$ws.Cells.Item(500, 15).Value2 = 0
$ws.Cells.Item(500, 16).Value2 = 0
$ws.Cells.Item(500, 17).Value2 = 0

$ws.Cells.Item(500, 18).WrapText = $true
$ws.Cells.Item(500, 18).Value2 = ""
